# Auto-generated edit script: updates crypto price/volume table cells
# per commit 'Updated cryptos list on Sat Jul  6 07:27:18 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.586.02"
$ws.Range("E2").Value = "  +4.14%  "
# Row 3
$ws.Range("D3").Value = "3.010.51"
$ws.Range("E3").Value = "  +4.87%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.19%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.90%  "
# Row 7
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.04%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +14.89%  "
# Row 10
$ws.Range("E10").Value = "  +13.67%  "
# Row 11
$ws.Range("E11").Value = "  +7.04%  "
# Row 12
$ws.Range("E12").Value = "  +4.60%  "
# Row 13
$ws.Range("D13").Value = "3.526.24"
$ws.Range("E13").Value = "  +4.92%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.90%  "
# Row 15
$ws.Range("E15").Value = "  +15.16%  "
# Row 16
$ws.Range("D16").Value = "56.679.50"
$ws.Range("E16").Value = "  +4.30%  "
# Row 17
$ws.Range("D17").Value = "3.008.74"
$ws.Range("E17").Value = "  +5.04%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.03%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.68%  "
# Row 20
$ws.Range("E20").Value = "  +12.04%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.28%  "
# Row 22
$ws.Range("E22").Value = "  -0.25%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.478"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.80%  "
# Row 24
$ws.Range("E24").Value = "  +6.39%  "
# Row 25
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.37%  "
# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.34%  "
# Row 27
$ws.Range("D27").Value = "0.0₃0920"
$ws.Range("E27").Value = "  +14.31%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.39%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.67%  "
# Row 30
$ws.Range("E30").Value = "  +10.14%  "
# Row 31
$ws.Range("E31").Value = "  +10.05%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.64%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.50%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.31%  "
# Row 35
$ws.Range("E35").Value = "  +3.74%  "
# Row 36
$ws.Range("E36").Value = "  +4.90%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0676"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.54%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.29%  "
# Row 39
$ws.Range("D39").Value = "3.049.03"
$ws.Range("E39").Value = "  +5.32%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.59%  "
# Row 41
$ws.Range("E41").Value = "  +0.01%  "
# Row 42
$ws.Range("E42").Value = "  +7.93%  "
# Row 43
$ws.Range("D43").Value = "2.268.91"
$ws.Range("E43").Value = "  +11.05%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.47%  "
# Row 45
$ws.Range("E45").Value = "  +7.52%  "
# Row 46
$ws.Range("E46").Value = "  +6.89%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +25.26%  "
# Row 48
$ws.Range("E48").Value = "  +11.94%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.54%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.32%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0876"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.02%  "
